$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.717.14'
$ws.Range("E2").Value = '  -1.64%  '

# Row 3
$ws.Range("D3").Value = '3.383.45'
$ws.Range("E3").Value = '  -2.09%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.23'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.64'
$ws.Range("E6").Value = '  -3.79%  '

# Row 8
$ws.Range("D8").Value = '3.384.01'
$ws.Range("E8").Value = '  -2.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.474'
$ws.Range("E9").Value = '  -0.38%  '

# Row 10
$ws.Range("E10").Value = '  -2.09%  '

# Row 11
$ws.Range("E11").Value = '  -2.20%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.399'
$ws.Range("E12").Value = '  +2.01%  '

# Row 13
$ws.Range("D13").Value = '3.963.13'
$ws.Range("E13").Value = '  -2.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.40'
$ws.Range("E14").Value = '  +1.37%  '

# Row 15
$ws.Range("E15").Value = '  +1.46%  '

# Row 16
$ws.Range("E16").Value = '  -2.21%  '

# Row 17
$ws.Range("D17").Value = '3.381.20'
$ws.Range("E17").Value = '  -2.43%  '

# Row 18
$ws.Range("D18").Value = '60.824.83'
$ws.Range("E18").Value = '  -1.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  +0.17%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.10'
$ws.Range("E20").Value = '  -2.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.03'
$ws.Range("E21").Value = '  -5.83%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.32'
$ws.Range("E22").Value = '  -1.39%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.562'
$ws.Range("E23").Value = '  -0.88%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.65'
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("E25").Value = '  +0.07%  '

# Row 26
$ws.Range("E26").Value = '  -5.21%  '

# Row 27
$ws.Range("D27").Value = '3.523.33'
$ws.Range("E27").Value = '  -2.09%  '

# Row 28
$ws.Range("E28").Value = '  -2.55%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.44'
$ws.Range("E29").Value = '  -3.41%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.35%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.02'
$ws.Range("E31").Value = '  -2.44%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.14'
$ws.Range("E32").Value = '  -2.52%  '

# Row 33
$ws.Range("E33").Value = '  -3.68%  '

# Row 35
$ws.Range("E35").Value = '  -2.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.99'
$ws.Range("E36").Value = '  -0.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.73'
$ws.Range("E37").Value = '  -0.25%  '

# Row 38
$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.415.67'
$ws.Range("E38").Value = '  -1.89%  '

# Row 39
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.01'
$ws.Range("E39").Value = '  -2.89%  '

# Row 40
$ws.Range("E40").Value = '  -4.80%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.91'
$ws.Range("E41").Value = '  +1.77%  '

# Row 42
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0777'
$ws.Range("E42").Value = '  -1.02%  '

# Row 43
$ws.Range("E43").Value = '  -3.05%  '

# Row 44
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.43'
$ws.Range("E45").Value = '  -1.98%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.77'
$ws.Range("E46").Value = '  -1.93%  '

# Row 47
$ws.Range("E47").Value = '  -2.94%  '

# Row 48
$ws.Range("D48").Value = '2.520.51'
$ws.Range("E48").Value = '  -2.29%  '

# Row 49
$ws.Range("E49").Value = '  -3.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.74'
$ws.Range("E50").Value = '  +3.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.86'
$ws.Range("E51").Value = '  -1.25%  '
